# Apply updated crypto prices / 1h volume percentages (cryptos.xlsx refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.027.71'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '1.591.64'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = "'" + '210.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = "'" + '0.481'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.78%  '
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").Value = "'" + '17.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.44%  '
$ws.Range("E11").Value = '  +2.33%  '
$ws.Range("D12").Value = '1.814.24'
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").Value = '1.593.39'
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '26.022.11'
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").Value = "'" + '60.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '0.0₃0723'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").Value = "'" + '201.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.55%  '
$ws.Range("D21").Value = "'" + '4.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.71%  '
$ws.Range("E22").Value = '  -1.30%  '
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("D24").Value = "'" + '1.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +16.37%  '
$ws.Range("D25").Value = "'" + '143.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  -7.88%  '
$ws.Range("D28").Value = "'" + '15.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("E29").Value = '  +0.64%  '
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("D31").Value = "'" + '0.0474'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("E33").Value = '  -2.62%  '
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("D36").Value = '1.128.98'
$ws.Range("E36").Value = '  +2.92%  '
$ws.Range("E37").Value = '  +8.27%  '
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").Value = "'" + '0.791'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.09%  '
$ws.Range("D40").Value = "'" + '2.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("E41").Value = '  -1.82%  '
$ws.Range("E42").Value = '  -4.00%  '
$ws.Range("D43").Value = "'" + '5.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("D44").Value = '1.724.62'
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("D45").Value = "'" + '92.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.54%  '
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("D47").Value = "'" + '53.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.16%  '
$ws.Range("D48").Value = "'" + '0.0502'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.09%  '
$ws.Range("E49").Value = '  -0.33%  '
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = '0.0₇0923'
$ws.Range("E51").Value = '  -17.25%  '
